{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph directly above it) that the site generator\n// used to append after the bibliography entries.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"text\");\nawait context.sync();\n\nconst items = paragraphs.items;\n\n// Locate the two footer paragraphs by their text content.\nlet jupiterIdx = -1;\nlet copyrightIdx = -1;\nfor (let i = 0; i < items.length; i++) {\n  const t = items[i].text;\n  if (jupiterIdx === -1 && t.indexOf(\"Ver no Jupiter\") !== -1) {\n    jupiterIdx = i;\n  }\n  if (copyrightIdx === -1 && t.indexOf(\"Powered by Jekyll\") !== -1) {\n    copyrightIdx = i;\n  }\n}\n\nif (jupiterIdx === -1 || copyrightIdx === -1) {\n  throw new Error(\"Could not locate the footer paragraphs to remove.\");\n}\n\n// The blank paragraph that immediately precedes \"Ver no Jupiter ...\" is\n// also removed (it only separated the bibliography from the footer).\nlet blankIdx = jupiterIdx - 1;\nif (blankIdx >= 0 && items[blankIdx].text !== \"\") {\n  blankIdx = -1;\n}\n\nconst toDelete = [];\nif (blankIdx !== -1) toDelete.push(items[blankIdx]);\ntoDelete.push(items[jupiterIdx]);\ntoDelete.push(items[copyrightIdx]);\n\n// Delete from the end backwards so earlier indices stay valid.\nfor (let i = toDelete.length - 1; i >= 0; i--) {\n  toDelete[i].delete();\n}\n\nawait context.sync();\n", "ps1": "# Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n# (and the blank paragraph directly above it) that the site generator\n# used to append after the bibliography entries.\n\n$d = $word.ActiveDocument\n\n$jupiterIdx = -1\n$copyrightIdx = -1\n\n$count = $d.Paragraphs.Count\nfor ($i = 1; $i -le $count; $i++) {\n    $t = $d.Paragraphs.Item($i).Range.Text\n    if ($jupiterIdx -eq -1 -and $t -like \"*Ver no Jupiter*\") {\n        $jupiterIdx = $i\n    }\n    if ($copyrightIdx -eq -1 -and $t -like \"*Powered by Jekyll*\") {\n        $copyrightIdx = $i\n    }\n}\n\nif ($jupiterIdx -eq -1 -or $copyrightIdx -eq -1) {\n    throw \"Could not locate the footer paragraphs to remove.\"\n}\n\n# The blank paragraph immediately preceding \"Ver no Jupiter ...\" is also\n# removed (it only separated the bibliography from the footer).\n$blankIdx = -1\n$prev = $jupiterIdx - 1\nif ($prev -ge 1) {\n    $prevText = $d.Paragraphs.Item($prev).Range.Text.Trim()\n    if ($prevText -eq \"\") {\n        $blankIdx = $prev\n    }\n}\n\n# Build a single contiguous range spanning from the blank paragraph (or the\n# Jupiter paragraph if there is none) through the end of the copyright\n# paragraph, and delete it in one go.\n$startIdx = $jupiterIdx\nif ($blankIdx -ne -1) {\n    $startIdx = $blankIdx\n}\n\n$startRange = $d.Paragraphs.Item($startIdx).Range\n$endRange = $d.Paragraphs.Item($copyrightIdx).Range\n\n$deleteRange = $d.Range($startRange.Start, $endRange.End)\n$deleteRange.Delete()\n\n$d.Save()\n"}
